# Edit script: applies the bitacora.docx diff via Word COM-interop
$d = $word.ActiveDocument

# --- Change 1: merge the three runs making up "Viernes 29 de septiembre:"
# into a single run. The visible text does not change, so we first set it
# to a transient placeholder (forcing the engine to collapse/merge the
# matched runs into one), then rewrite it back to the final text.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Viernes 29 de septiembre:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Text = "Viernes 29 de septiembre:{PLACEHOLDER}"
    $rng2 = $d.Content
    $rng2.Find.Execute("Viernes 29 de septiembre:{PLACEHOLDER}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng2.Text = "Viernes 29 de septiembre:"
}

# --- Change 2: append the new "Viernes 13 de octubre" section (one blank
# divider paragraph, the date heading, and seven bulleted to-do items) at
# the end of the document, after "Ajuste de hiperparámetros".
#
# InsertXML on a range that sits at the very end of the document body
# clobbers the trailing paragraph's own content, so we first grow the
# document with two throwaway empty paragraphs (making our insertion
# point no longer "last in story"), inject the new content there, and
# finally delete the left-over helper paragraph it pushed to the end.
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$helper = $d.Paragraphs.Last
$hr = $helper.Range
$hr.Collapse(0)
$hr.InsertParagraphAfter()

$count = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($count - 1)
$tr = $targetPara.Range

$xmlFrag = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Viernes 13 de octubre:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Agregar etiquetas a los </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>cluster</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> y l</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>o</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>s gráficos de tendencia según los días o semana</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Agregar regularización</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Random search, bayesiana search</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Justificar muy bien la descomposición estacional</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Agregar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>label</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> a los eje</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> de todo</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">Agregar cuantas </w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">semanas </w:t></w:r><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>son el conjunto de prueba</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:lastRenderedPageBreak/><w:t>Agregar cada cuanto se necesita reentrenar el modelo</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$tr.InsertXML($xmlFrag)

# Remove the trailing helper paragraph left after the injected content.
$trailing = $d.Paragraphs.Last
$trailing.Range.Delete()

Write-Output ("Done. ParagraphCount=" + $d.Paragraphs.Count)
